# Update "Unidades Pedido" (L) and "Diferencia Stock" (M) for the rows
# whose order quantities were recalculated so the summary sections all show up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 5;  L = 3;  M = 0 },
    @{ Row = 10; L = 1;  M = 0 },
    @{ Row = 15; L = 4;  M = 0 },
    @{ Row = 23; L = 10; M = 0 },
    @{ Row = 25; L = 4;  M = 0 },
    @{ Row = 26; L = 1;  M = 0 },
    @{ Row = 29; L = 2;  M = 0 },
    @{ Row = 30; L = 6;  M = 0 },
    @{ Row = 34; L = 1;  M = 0 },
    @{ Row = 35; L = 4;  M = 0 },
    @{ Row = 36; L = 1;  M = 0 },
    @{ Row = 40; L = 1;  M = 0 },
    @{ Row = 47; L = 3;  M = 0 },
    @{ Row = 48; L = 1;  M = 0 },
    @{ Row = 70; L = 1;  M = 0 },
    @{ Row = 74; L = 6;  M = 0 },
    @{ Row = 81; L = 9;  M = 0 },
    @{ Row = 82; L = 6;  M = 0 },
    @{ Row = 85; L = 7;  M = 0 }
)

foreach ($u in $updates) {
    $ws.Range("L$($u.Row)").Value = $u.L
    $ws.Range("M$($u.Row)").Value = $u.M
}

# Update the summary metrics (resumen_pedido)
$ws.Range("C88").Value = 214
$ws.Range("C99").Value = 0
